$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row pairs whose data (columns B through AC) need to be swapped between them.
# Column A (the sequential id/index column) stays untouched.
$pairs = @(
    @(9, 10),
    @(49, 50),
    @(99, 100),
    @(122, 123)
)

$firstCol = 2   # Column B
$lastCol = 29   # Column AC

foreach ($pair in $pairs) {
    $r1 = $pair[0]
    $r2 = $pair[1]

    $range1 = $ws.Range($ws.Cells.Item($r1, $firstCol), $ws.Cells.Item($r1, $lastCol))
    $range2 = $ws.Range($ws.Cells.Item($r2, $firstCol), $ws.Cells.Item($r2, $lastCol))

    $values1 = $range1.Value()
    $values2 = $range2.Value()

    $range1.Value = $values2
    $range2.Value = $values1
}
